$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 2981.7
$ws.Range("I76").Value = 2976.258
$ws.Range("J76").Value = 3054
$ws.Range("K76").Value = 2976.258
$ws.Range("L76").Value = 3054
$ws.Range("M76").Value = -2661.258
$ws.Range("N76").Value = -3684
$ws.Range("H79").Value = 2981.7
$ws.Range("I79").Value = 2976.258
$ws.Range("J79").Value = 3054
$ws.Range("K79").Value = 2976.258
$ws.Range("L79").Value = 3054
$ws.Range("M79").Value = -1884.258
$ws.Range("N79").Value = -5238
$ws.Range("H98").Value = 1172.1072
$ws.Range("I98").Value = 838.96155
$ws.Range("K98").Value = 838.96155
$ws.Range("M98").Value = 659.03845
$ws.Range("H122").Value = 1172.1072
$ws.Range("I122").Value = 838.96155
$ws.Range("K122").Value = 2516.88465
$ws.Range("M122").Value = -66.88464999999997
$ws.Range("H128").Value = 21666.666
$ws.Range("J128").Value = 21666.666
$ws.Range("L128").Value = 21666.666
$ws.Range("N128").Value = -31626.666
$ws.Range("H129").Value = 966.45
$ws.Range("I129").Value = 782
$ws.Range("J129").Value = 990.81134
$ws.Range("K129").Value = 2346
$ws.Range("L129").Value = 2972.43402
$ws.Range("M129").Value = 2654
$ws.Range("N129").Value = -12972.43402
$ws.Range("H131").Value = 2598.3333
$ws.Range("I131").Value = 2085
$ws.Range("J131").Value = 3625
$ws.Range("K131").Value = 6255
$ws.Range("L131").Value = 10875
$ws.Range("M131").Value = -1215
$ws.Range("N131").Value = -20955
$ws.Range("H132").Value = 1022.0417
$ws.Range("I132").Value = 614.7059
$ws.Range("J132").Value = 2011.2858
$ws.Range("K132").Value = 1844.1177
$ws.Range("L132").Value = 6033.857400000001
$ws.Range("M132").Value = 685.8822999999998
$ws.Range("N132").Value = -11093.8574
$ws.Range("H133").Value = 40780
$ws.Range("J133").Value = 40780
$ws.Range("L133").Value = 40780
$ws.Range("N133").Value = -50900
$ws.Range("H137").Value = 1428.9272
$ws.Range("I137").Value = 1192.375
$ws.Range("J137").Value = 2059.7334
$ws.Range("K137").Value = 3577.125
$ws.Range("L137").Value = 6179.2002
$ws.Range("M137").Value = -1027.125
$ws.Range("N137").Value = -11279.2002
$ws.Range("H141").Value = 971.5965
$ws.Range("I141").Value = 734.72546
$ws.Range("J141").Value = 2985
$ws.Range("K141").Value = 2204.17638
$ws.Range("L141").Value = 8955
$ws.Range("M141").Value = 2975.82362
$ws.Range("N141").Value = -19315

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1163.0588
$ws.Range("I2").Value = 923.25
$ws.Range("J2").Value = 5000
$ws.Range("K2").Value = 923.25
$ws.Range("L2").Value = 5000
$ws.Range("M2").Value = -810.25
$ws.Range("N2").Value = -5226
$ws.Range("H31").Value = 20200
$ws.Range("I31").Value = 14500
$ws.Range("J31").Value = 43000
$ws.Range("K31").Value = 14500
$ws.Range("L31").Value = 43000
$ws.Range("M31").Value = -14206
$ws.Range("N31").Value = -43588
$ws.Range("H61").Value = 1394.9131
$ws.Range("I61").Value = 1268.8485
$ws.Range("J61").Value = 1714.9231
$ws.Range("K61").Value = 1268.8485
$ws.Range("L61").Value = 1714.9231
$ws.Range("M61").Value = -1056.8485
$ws.Range("N61").Value = -2138.9231
$ws.Range("H74").Value = 1044.8372
$ws.Range("I74").Value = 993.2963
$ws.Range("J74").Value = 1131.8125
$ws.Range("K74").Value = 993.2963
$ws.Range("L74").Value = 1131.8125
$ws.Range("M74").Value = -119.2963
$ws.Range("N74").Value = -2879.8125
$ws.Range("H77").Value = 1044.8372
$ws.Range("I77").Value = 993.2963
$ws.Range("J77").Value = 1131.8125
$ws.Range("K77").Value = 4966.4815
$ws.Range("L77").Value = 5659.0625
$ws.Range("M77").Value = -598.4814999999999
$ws.Range("N77").Value = -14395.0625
$ws.Range("H114").Value = 34999
$ws.Range("J114").Value = 34999
$ws.Range("L114").Value = 34999
$ws.Range("N114").Value = -43677
$ws.Range("H116").Value = 1163.0588
$ws.Range("I116").Value = 923.25
$ws.Range("J116").Value = 5000
$ws.Range("K116").Value = 923.25
$ws.Range("L116").Value = 5000
$ws.Range("M116").Value = 1370.75
$ws.Range("N116").Value = -9588
$ws.Range("H132").Value = 1494591.5
$ws.Range("I132").Value = 1403.1296
$ws.Range("J132").Value = 7697066
$ws.Range("K132").Value = 4209.3888
$ws.Range("L132").Value = 23091198
$ws.Range("M132").Value = -1679.3888
$ws.Range("N132").Value = -23096258
$ws.Range("H136").Value = 1394.9131
$ws.Range("I136").Value = 1268.8485
$ws.Range("J136").Value = 1714.9231
$ws.Range("K136").Value = 3806.5455
$ws.Range("L136").Value = 5144.7693
$ws.Range("M136").Value = -1256.5455
$ws.Range("N136").Value = -10244.7693

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1163.0588
$ws.Range("I3").Value = 923.25
$ws.Range("J3").Value = 5000
$ws.Range("K3").Value = 923.25
$ws.Range("L3").Value = 5000
$ws.Range("M3").Value = -809.25
$ws.Range("N3").Value = -5228
$ws.Range("H134").Value = 1439.5074
$ws.Range("I134").Value = 1040.2916
$ws.Range("J134").Value = 2448.0527
$ws.Range("K134").Value = 3120.8748
$ws.Range("L134").Value = 7344.158100000001
$ws.Range("M134").Value = -585.8748000000001
$ws.Range("N134").Value = -12414.1581

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 369.46667
$ws.Range("I22").Value = 325.1
$ws.Range("K22").Value = 325.1
$ws.Range("M22").Value = 24.89999999999998
$ws.Range("H31").Value = 5246.7236
$ws.Range("I31").Value = 1260.5686
$ws.Range("J31").Value = 13378.48
$ws.Range("K31").Value = 1260.5686
$ws.Range("L31").Value = 13378.48
$ws.Range("M31").Value = -965.5686000000001
$ws.Range("N31").Value = -13968.48
$ws.Range("H34").Value = 5246.7236
$ws.Range("I34").Value = 1260.5686
$ws.Range("J34").Value = 13378.48
$ws.Range("K34").Value = 1260.5686
$ws.Range("L34").Value = 13378.48
$ws.Range("M34").Value = -1058.5686
$ws.Range("N34").Value = -13782.48
$ws.Range("H99").Value = 5689954
$ws.Range("I99").Value = 7848.0625
$ws.Range("J99").Value = 20842236
$ws.Range("K99").Value = 7848.0625
$ws.Range("L99").Value = 20842236
$ws.Range("M99").Value = -6350.0625
$ws.Range("N99").Value = -20845232
$ws.Range("H126").Value = 5689954
$ws.Range("I126").Value = 7848.0625
$ws.Range("J126").Value = 20842236
$ws.Range("K126").Value = 23544.1875
$ws.Range("L126").Value = 62526708
$ws.Range("M126").Value = -21074.1875
$ws.Range("N126").Value = -62531648
$ws.Range("H132").Value = 1687.6364
$ws.Range("I132").Value = 1294.5122
$ws.Range("J132").Value = 2838.9285
$ws.Range("K132").Value = 3883.536599999999
$ws.Range("L132").Value = 8516.7855
$ws.Range("M132").Value = -1353.536599999999
$ws.Range("N132").Value = -13576.7855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 172974.83
$ws.Range("J5").Value = 417687.03
$ws.Range("L5").Value = 1253061.09
$ws.Range("N5").Value = -1253285.09
$ws.Range("H7").Value = 99.71429000000001
$ws.Range("I7").Value = 99.59999999999999
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 298.8
$ws.Range("L7").Value = 300
$ws.Range("M7").Value = -186.8
$ws.Range("N7").Value = -524
$ws.Range("H135").Value = 172974.83
$ws.Range("J135").Value = 417687.03
$ws.Range("L135").Value = 3759183.27
$ws.Range("N135").Value = -3764253.27
$ws.Range("H140").Value = 35716276
$ws.Range("I140").Value = 35716276
$ws.Range("K140").Value = 107148828
$ws.Range("M140").Value = -107143648
$ws.Range("H141").Value = 13148.546
$ws.Range("I141").Value = 15384.429
$ws.Range("J141").Value = 9235.75
$ws.Range("K141").Value = 46153.287
$ws.Range("L141").Value = 27707.25
$ws.Range("M141").Value = -40973.287
$ws.Range("N141").Value = -38067.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2046.9
$ws.Range("I102").Value = 1793.8
$ws.Range("J102").Value = 2300
$ws.Range("K102").Value = 1793.8
$ws.Range("L102").Value = 2300
$ws.Range("M102").Value = -171.8
$ws.Range("N102").Value = -5544
$ws.Range("H126").Value = 5169.2
$ws.Range("I126").Value = 8526.4
$ws.Range("J126").Value = 2651.3
$ws.Range("K126").Value = 25579.2
$ws.Range("L126").Value = 7953.900000000001
$ws.Range("M126").Value = -23109.2
$ws.Range("N126").Value = -12893.9
$ws.Range("H132").Value = 2144.4634
$ws.Range("I132").Value = 1831.2963
$ws.Range("J132").Value = 2748.4285
$ws.Range("K132").Value = 5493.8889
$ws.Range("L132").Value = 8245.2855
$ws.Range("M132").Value = -2963.8889
$ws.Range("N132").Value = -13305.2855
$ws.Range("H133").Value = 40290
$ws.Range("J133").Value = 40290
$ws.Range("L133").Value = 40290
$ws.Range("N133").Value = -50410

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 41942.24
$ws.Range("I7").Value = 51782.6
$ws.Range("K7").Value = 51782.6
$ws.Range("M7").Value = -51670.6
$ws.Range("H40").Value = 41670476
$ws.Range("I40").Value = 55558130
$ws.Range("K40").Value = 55558130
$ws.Range("M40").Value = -55557994
$ws.Range("H122").Value = 1491530.8
$ws.Range("I122").Value = 1934177.8
$ws.Range("J122").Value = 2627.0908
$ws.Range("K122").Value = 5802533.4
$ws.Range("L122").Value = 7881.2724
$ws.Range("M122").Value = -5800083.4
$ws.Range("N122").Value = -12781.2724
$ws.Range("H126").Value = 41942.24
$ws.Range("I126").Value = 51782.6
$ws.Range("K126").Value = 155347.8
$ws.Range("M126").Value = -152877.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1918.5807
$ws.Range("I122").Value = 1880.25
$ws.Range("J122").Value = 2050
$ws.Range("K122").Value = 5640.75
$ws.Range("L122").Value = 6150
$ws.Range("M122").Value = -3190.75
$ws.Range("N122").Value = -11050
$ws.Range("H136").Value = 5557350.5
$ws.Range("I136").Value = 1729.6719
$ws.Range("J136").Value = 19232724
$ws.Range("K136").Value = 5189.0157
$ws.Range("L136").Value = 57698172
$ws.Range("M136").Value = -2639.0157
$ws.Range("N136").Value = -57703272
